$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Set tab colors on all three sheets (rgb 00FFFFFF -> white, alpha 00)
$ws1.Tab.Color = 16777215
$ws2.Tab.Color = 16777215
$ws3.Tab.Color = 16777215

# Update formulas on Sheet1 row 3 and row 4
$ws1.Range("B3").Formula = "=PI()/6"
$ws1.Range("F3").Formula = "=SQRT(30*30+62*62)"
$ws1.Range("B4").Formula = "=5*PI()/6"

# Fix selections (active cell) on each sheet
$ws1.Range("B5").Select()
$ws2.Range("A1").Select()
$ws3.Range("A1").Select()

# Make Sheet1 active/selected at the end (tabSelected=true in sheet1)
$ws1.Activate()
$ws1.Range("B5").Select()
